# Apply updated "dSF" (column F) values for specific rows as part of a
# data repull / recalculation pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 1
    11 = -1
    20 = 0
    23 = -2
    35 = -4
    37 = 2
    50 = 0
    53 = -2
    55 = 0
    58 = -2
    66 = 2
    71 = -1
    80 = 6
    82 = -2
    83 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
